$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.86382305527799
$ws.Cells.Item(2, 3).Value = 11.56520776152982
$ws.Cells.Item(2, 4).Value = 5.979087027223589
$ws.Cells.Item(2, 5).Value = 12.57996649746447
$ws.Cells.Item(2, 7).Value = 27.16832293086387
$ws.Cells.Item(2, 8).Value = 13.88820697996222
$ws.Cells.Item(2, 9).Value = 21.48495719688379
$ws.Cells.Item(2, 12).Value = 10.0083676588922
$ws.Cells.Item(2, 13).Value = 14.32865406225001
$ws.Cells.Item(2, 15).Value = 20.92838059598559
$ws.Cells.Item(3, 2).Value = 14.24410338228091
$ws.Cells.Item(3, 3).Value = 11.31158168539129
$ws.Cells.Item(3, 4).Value = 5.859132931411608
$ws.Cells.Item(3, 5).Value = 12.63219431836582
$ws.Cells.Item(3, 7).Value = 27.24827256104453
$ws.Cells.Item(3, 8).Value = 13.94667884802119
$ws.Cells.Item(3, 9).Value = 21.62763977056519
$ws.Cells.Item(3, 12).Value = 10.01707382128698
$ws.Cells.Item(3, 13).Value = 14.18703841361733
$ws.Cells.Item(3, 15).Value = 21.02103435551991
$ws.Cells.Item(4, 2).Value = 13.84980895604515
$ws.Cells.Item(4, 3).Value = 11.15228360332542
$ws.Cells.Item(4, 4).Value = 5.785990792824317
$ws.Cells.Item(4, 5).Value = 12.66595129450744
$ws.Cells.Item(4, 7).Value = 27.30925428398834
$ws.Cells.Item(4, 8).Value = 13.98541053627015
$ws.Cells.Item(4, 9).Value = 21.72043853021644
$ws.Cells.Item(4, 12).Value = 10.02383823876836
$ws.Cells.Item(4, 13).Value = 14.10081033826655
$ws.Cells.Item(4, 15).Value = 21.08380136186442
$ws.Cells.Item(5, 2).Value = 13.68587244037144
$ws.Cells.Item(5, 3).Value = 11.08653630134441
$ws.Cells.Item(5, 4).Value = 5.756359565720156
$ws.Cells.Item(5, 5).Value = 12.68013347566257
$ws.Cells.Item(5, 7).Value = 27.33707600681361
$ws.Cells.Item(5, 8).Value = 14.00190481248795
$ws.Cells.Item(5, 9).Value = 21.75956034960362
$ws.Cells.Item(5, 12).Value = 10.02695200507387
$ws.Cells.Item(5, 13).Value = 14.06588427980608
$ws.Cells.Item(5, 15).Value = 21.11085195044968
$ws.Cells.Item(6, 2).Value = 13.65846071780506
$ws.Cells.Item(6, 3).Value = 11.0755707487081
$ws.Cells.Item(6, 4).Value = 5.751451248334364
$ws.Cells.Item(6, 5).Value = 12.68251417868212
$ws.Cells.Item(6, 7).Value = 27.34187463742103
$ws.Cells.Item(6, 8).Value = 14.00468658295514
$ws.Cells.Item(6, 9).Value = 21.76613535961803
$ws.Cells.Item(6, 12).Value = 10.02749062986211
$ws.Cells.Item(6, 13).Value = 14.06009855479556
$ws.Cells.Item(6, 15).Value = 21.11543245833855
$ws.Cells.Item(7, 2).Value = 13.84761094927877
$ws.Cells.Item(7, 3).Value = 11.15140019165575
$ws.Cells.Item(7, 4).Value = 5.785590406198023
$ws.Cells.Item(7, 5).Value = 12.66614083407898
$ws.Cells.Item(7, 7).Value = 27.30961749406305
$ws.Cells.Item(7, 8).Value = 13.9856301068582
$ws.Cells.Item(7, 9).Value = 21.7209608541937
$ws.Cells.Item(7, 12).Value = 10.02387878523995
$ws.Cells.Item(7, 13).Value = 14.10033841261657
$ws.Cells.Item(7, 15).Value = 21.08416022054014
$ws.Cells.Item(8, 2).Value = 14.65312101179831
$ws.Cells.Item(8, 3).Value = 11.47853700326698
$ws.Cells.Item(8, 4).Value = 5.937648685232857
$ws.Cells.Item(8, 5).Value = 12.59762488025979
$ws.Cells.Item(8, 7).Value = 27.19341028107136
$ws.Cells.Item(8, 8).Value = 13.90778031538083
$ws.Cells.Item(8, 9).Value = 21.53307686628265
$ws.Cells.Item(8, 12).Value = 10.01107535651261
$ws.Cells.Item(8, 13).Value = 14.27969295246636
$ws.Cells.Item(8, 15).Value = 20.95910461123263
$ws.Cells.Item(9, 2).Value = 16.11553629522517
$ws.Cells.Item(9, 3).Value = 12.08900767842801
$ws.Cells.Item(9, 4).Value = 6.237877307633465
$ws.Cells.Item(9, 5).Value = 12.47660830236519
$ws.Cells.Item(9, 7).Value = 27.06066062287129
$ws.Cells.Item(9, 8).Value = 13.7776035157371
$ws.Cells.Item(9, 9).Value = 21.20581551893623
$ws.Cells.Item(9, 12).Value = 9.997204081702511
$ws.Cells.Item(9, 13).Value = 14.63580828443741
$ws.Cells.Item(9, 15).Value = 20.760738503309
$ws.Cells.Item(10, 2).Value = 17.1097149052258
$ws.Cells.Item(10, 3).Value = 12.51522616462274
$ws.Cells.Item(10, 4).Value = 6.457100170105859
$ws.Cells.Item(10, 5).Value = 12.3957508357594
$ws.Cells.Item(10, 7).Value = 27.02204386977153
$ws.Cells.Item(10, 8).Value = 13.69571705936493
$ws.Cells.Item(10, 9).Value = 20.99046645205752
$ws.Cells.Item(10, 12).Value = 9.993831292148053
$ws.Cells.Item(10, 13).Value = 14.8983194868867
$ws.Cells.Item(10, 15).Value = 20.64388373671845
$ws.Cells.Item(11, 2).Value = 17.54311786945895
$ws.Cells.Item(11, 3).Value = 12.70364614533142
$ws.Cells.Item(11, 4).Value = 6.556018421101848
$ws.Cells.Item(11, 5).Value = 12.36069819903405
$ws.Cells.Item(11, 7).Value = 27.01741996282442
$ws.Cells.Item(11, 8).Value = 13.66146151960847
$ws.Cells.Item(11, 9).Value = 20.89794568240518
$ws.Cells.Item(11, 12).Value = 9.993769296153935
$ws.Cells.Item(11, 13).Value = 15.01754923186129
$ws.Cells.Item(11, 15).Value = 20.59705979915659
$ws.Cells.Item(12, 2).Value = 17.70442782804727
$ws.Cells.Item(12, 3).Value = 12.7741613850164
$ws.Cells.Item(12, 4).Value = 6.593319756678811
$ws.Cells.Item(12, 5).Value = 12.34767213055801
$ws.Cells.Item(12, 7).Value = 27.01753894126538
$ws.Cells.Item(12, 8).Value = 13.64892141788191
$ws.Cells.Item(12, 9).Value = 20.86369345554918
$ws.Cells.Item(12, 12).Value = 9.993956689001644
$ws.Cells.Item(12, 13).Value = 15.06263963465854
$ws.Cells.Item(12, 15).Value = 20.58024466631055
$ws.Cells.Item(13, 2).Value = 17.66981325092717
$ws.Cells.Item(13, 3).Value = 12.75901255543944
$ws.Cells.Item(13, 4).Value = 6.585293847448115
$ws.Cells.Item(13, 5).Value = 12.35046653330634
$ws.Cells.Item(13, 7).Value = 27.01743003746513
$ws.Cells.Item(13, 8).Value = 13.65160293523403
$ws.Cells.Item(13, 9).Value = 20.87103543133805
$ws.Cells.Item(13, 12).Value = 9.993906966481934
$ws.Cells.Item(13, 13).Value = 15.05293175640011
$ws.Cells.Item(13, 15).Value = 20.58382528135053
$ws.Cells.Item(14, 2).Value = 17.55644577421667
$ws.Cells.Item(14, 3).Value = 12.70946448405084
$ws.Cells.Item(14, 4).Value = 6.559090593395268
$ws.Cells.Item(14, 5).Value = 12.35962158069608
$ws.Cells.Item(14, 7).Value = 27.01739223328018
$ws.Cells.Item(14, 8).Value = 13.66042117941338
$ws.Cells.Item(14, 9).Value = 20.89511202631203
$ws.Cells.Item(14, 12).Value = 9.993780492080472
$ws.Cells.Item(14, 13).Value = 15.02126021834584
$ws.Cells.Item(14, 15).Value = 20.59565801899744
$ws.Cells.Item(15, 2).Value = 17.48663624973829
$ws.Cells.Item(15, 3).Value = 12.67900469069407
$ws.Cells.Item(15, 4).Value = 6.54301874552497
$ws.Cells.Item(15, 5).Value = 12.36526152319322
$ws.Cells.Item(15, 7).Value = 27.01761281771062
$ws.Cells.Item(15, 8).Value = 13.66587886343938
$ws.Cells.Item(15, 9).Value = 20.90996167920535
$ws.Cells.Item(15, 12).Value = 9.99373045775185
$ws.Cells.Item(15, 13).Value = 15.00185181250826
$ws.Cells.Item(15, 15).Value = 20.60302536589439
$ws.Cells.Item(16, 2).Value = 17.08100278232849
$ws.Cells.Item(16, 3).Value = 12.50279816010508
$ws.Cells.Item(16, 4).Value = 6.45061584710267
$ws.Cells.Item(16, 5).Value = 12.39807631867447
$ws.Cells.Item(16, 7).Value = 27.02260730500949
$ws.Cells.Item(16, 8).Value = 13.69801609636537
$ws.Cells.Item(16, 9).Value = 20.99662248170156
$ws.Cells.Item(16, 12).Value = 9.993864899809532
$ws.Cells.Item(16, 13).Value = 14.89052098990269
$ws.Cells.Item(16, 15).Value = 20.64707169988687
$ws.Cells.Item(17, 2).Value = 16.82725530891032
$ws.Cells.Item(17, 3).Value = 12.39326443223025
$ws.Cells.Item(17, 4).Value = 6.393693157953982
$ws.Cells.Item(17, 5).Value = 12.41864936745757
$ws.Cells.Item(17, 7).Value = 27.02899329621687
$ws.Cells.Item(17, 8).Value = 13.71849908838471
$ws.Cells.Item(17, 9).Value = 21.05118054851009
$ws.Cells.Item(17, 12).Value = 9.994323892163656
$ws.Cells.Item(17, 13).Value = 14.82215188264693
$ws.Cells.Item(17, 15).Value = 20.67571876692049
$ws.Cells.Item(18, 2).Value = 16.67953672993436
$ws.Cells.Item(18, 3).Value = 12.32975144395649
$ws.Cells.Item(18, 4).Value = 6.360879041062486
$ws.Cells.Item(18, 5).Value = 12.43064532854422
$ws.Cells.Item(18, 7).Value = 27.0338842402301
$ws.Cells.Item(18, 8).Value = 13.73056220613667
$ws.Cells.Item(18, 9).Value = 21.08307313731129
$ws.Cells.Item(18, 12).Value = 9.99472652839569
$ws.Cells.Item(18, 13).Value = 14.78281199849002
$ws.Cells.Item(18, 15).Value = 20.69279167299235
$ws.Cells.Item(19, 2).Value = 16.62922107963503
$ws.Cells.Item(19, 3).Value = 12.30816063864836
$ws.Cells.Item(19, 4).Value = 6.349757367856512
$ws.Cells.Item(19, 5).Value = 12.43473496593666
$ws.Cells.Item(19, 7).Value = 27.03574909669079
$ws.Cells.Item(19, 8).Value = 13.73469494541829
$ws.Cells.Item(19, 9).Value = 21.09395937553203
$ws.Cells.Item(19, 12).Value = 9.994886689385426
$ws.Cells.Item(19, 13).Value = 14.76949047077357
$ws.Cells.Item(19, 15).Value = 20.69867444342655
$ws.Cells.Item(20, 2).Value = 16.85445108969325
$ws.Cells.Item(20, 3).Value = 12.40497784968615
$ws.Cells.Item(20, 4).Value = 6.399760604805592
$ws.Cells.Item(20, 5).Value = 12.41644248135278
$ws.Cells.Item(20, 7).Value = 27.02818737813704
$ws.Cells.Item(20, 8).Value = 13.71628946168846
$ws.Cells.Item(20, 9).Value = 21.04531973007209
$ws.Cells.Item(20, 12).Value = 9.99426068920708
$ws.Cells.Item(20, 13).Value = 14.82943175965134
$ws.Cells.Item(20, 15).Value = 20.67260752854838
$ws.Cells.Item(21, 2).Value = 17.58982153694969
$ws.Cells.Item(21, 3).Value = 12.72404099471121
$ws.Cells.Item(21, 4).Value = 6.56679169480859
$ws.Cells.Item(21, 5).Value = 12.35692581051017
$ws.Cells.Item(21, 7).Value = 27.01735252982189
$ws.Cells.Item(21, 8).Value = 13.65781932200749
$ws.Cells.Item(21, 9).Value = 20.8880188795294
$ws.Cells.Item(21, 12).Value = 9.993811924897473
$ws.Cells.Item(21, 13).Value = 15.03056478141059
$ws.Cells.Item(21, 15).Value = 20.59215755399166
$ws.Cells.Item(22, 2).Value = 18.05401218261721
$ws.Cells.Item(22, 3).Value = 12.92767832264986
$ws.Cells.Item(22, 4).Value = 6.675024620181423
$ws.Cells.Item(22, 5).Value = 12.31947092923945
$ws.Cells.Item(22, 7).Value = 27.02117313856916
$ws.Cells.Item(22, 8).Value = 13.62212263366064
$ws.Cells.Item(22, 9).Value = 20.78978024996395
$ws.Cells.Item(22, 12).Value = 9.994747288433189
$ws.Cells.Item(22, 13).Value = 15.16165651403867
$ws.Cells.Item(22, 15).Value = 20.54492074840686
$ws.Cells.Item(23, 2).Value = 17.80779592250748
$ws.Cells.Item(23, 3).Value = 12.8194557532185
$ws.Cells.Item(23, 4).Value = 6.617356772000178
$ws.Cells.Item(23, 5).Value = 12.33932967203756
$ws.Cells.Item(23, 7).Value = 27.01813424170082
$ws.Cells.Item(23, 8).Value = 13.64094398201899
$ws.Cells.Item(23, 9).Value = 20.84179393865952
$ws.Cells.Item(23, 12).Value = 9.994135946087431
$ws.Cells.Item(23, 13).Value = 15.09173389076072
$ws.Cells.Item(23, 15).Value = 20.56964144333812
$ws.Cells.Item(24, 2).Value = 16.84216158756414
$ws.Cells.Item(24, 3).Value = 12.39968389374963
$ws.Cells.Item(24, 4).Value = 6.397017784646511
$ws.Cells.Item(24, 5).Value = 12.41743969041559
$ws.Cells.Item(24, 7).Value = 27.02854793556544
$ws.Cells.Item(24, 8).Value = 13.71728753945575
$ws.Cells.Item(24, 9).Value = 21.0479677663187
$ws.Cells.Item(24, 12).Value = 9.994288830992563
$ws.Cells.Item(24, 13).Value = 14.82614062935999
$ws.Cells.Item(24, 15).Value = 20.67401224047697
$ws.Cells.Item(25, 2).Value = 15.73345449636791
$ws.Cells.Item(25, 3).Value = 11.92754522849003
$ws.Cells.Item(25, 4).Value = 6.156713297433891
$ws.Cells.Item(25, 5).Value = 12.50792642438084
$ws.Cells.Item(25, 7).Value = 27.08628510623515
$ws.Cells.Item(25, 8).Value = 13.8104077295232
$ws.Cells.Item(25, 9).Value = 21.28994174822689
$ws.Cells.Item(25, 12).Value = 9.999756655448856
$ws.Cells.Item(25, 13).Value = 14.53918535950729
$ws.Cells.Item(25, 15).Value = 20.80935142295522
